# Fix the velocity-profile formula equation (2nd OMath in the document):
#   ... A^2/(r2^2-r1^2) [ q/2pi + 1 - 1/(r2^2-r1^2) ] (...)
# becomes
#   ... A^2/(r2^2-r1^2) ( q/2pi + (r2^2-r1^2)/2 ) - 1/(r2^2-r1^2) (...)
# Also strips the redundant explicit <w:sz>/<w:szCs val="24"/> run
# properties from that whole equation (matching the author's re-save
# of the formula).

$d = $word.ActiveDocument
$om = $d.OMaths.Item(2)
$rng = $om.Range
$xml = $rng.WordOpenXML

$startTag = "<m:oMathPara>"
$endTag = "</m:oMathPara>"
$startIdx = $xml.IndexOf($startTag)
$endIdx = $xml.IndexOf($endTag) + $endTag.Length
$head = $xml.Substring(0, $startIdx)
$mathBlock = $xml.Substring($startIdx, $endIdx - $startIdx)
$tail = $xml.Substring($endIdx)

# Strip the redundant explicit run-size overrides (sz/szCs = 24
# half-points) throughout this equation -- the author's re-save
# dropped every one of them in this formula.
$mathBlock = $mathBlock.Replace('<w:sz w:val="24"/>', '')
$mathBlock = $mathBlock.Replace('<w:szCs w:val="24"/>', '')

$rPr  = '<w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr>'
$rPrI = '<w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:i/></w:rPr>'
$fPr  = '<m:fPr><m:ctrlPr>' + $rPrI + '</m:ctrlPr></m:fPr>'
$sSubSupPr = '<m:sSubSupPr><m:ctrlPr>' + $rPrI + '</m:ctrlPr></m:sSubSupPr>'

function Run([string]$t) {
    return '<m:r>' + $rPr + '<m:t>' + $t + '</m:t></m:r>'
}
function SubSup([string]$base, [string]$sub, [string]$sup) {
    return '<m:sSubSup>' + $sSubSupPr + '<m:e>' + (Run $base) + '</m:e>' +
        '<m:sub>' + (Run $sub) + '</m:sub><m:sup>' + (Run $sup) + '</m:sup></m:sSubSup>'
}

# Fraction q / 2*pi (kept identical to the original)
$fracQ2Pi = '<m:f>' + $fPr + '<m:num>' + (Run "q") + '</m:num><m:den>' + (Run "2π") + '</m:den></m:f>'

# New fraction (r2^2 - r1^2) / 2
$numDiff = (SubSup "r" "2" "2") + (Run "-") + (SubSup "r" "1" "2")
$fracNew = '<m:f>' + $fPr + '<m:num>' + $numDiff + '</m:num><m:den>' + (Run "2") + '</m:den></m:f>'

# Fraction 1 / (r2^2 - r1^2) (kept identical to the original)
$denDiff = (SubSup "r" "2" "2") + (Run "-") + (SubSup "r" "1" "2")
$frac1Over = '<m:f>' + $fPr + '<m:num>' + (Run "1") + '</m:num><m:den>' + $denDiff + '</m:den></m:f>'

$replacement = (Run "(") + $fracQ2Pi + (Run "+") + $fracNew + (Run ")") + (Run "-") + $frac1Over

$oldSegStart = '<m:f><m:fPr><m:ctrlPr><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/><w:i/></w:rPr></m:ctrlPr></m:fPr><m:num><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:hAnsi="Cambria Math"/></w:rPr><m:t>q</m:t>'
$oldSegEndMarker = '<m:d><m:dPr><m:begChr m:val="["/>'

$segStartIdx = $mathBlock.IndexOf($oldSegStart)
$segEndIdx = $mathBlock.IndexOf($oldSegEndMarker)

if ($segStartIdx -lt 0 -or $segEndIdx -lt 0) {
    throw "Could not locate the velocity-profile fraction segment to patch."
}

$newMathBlock = $mathBlock.Substring(0, $segStartIdx) + $replacement + $mathBlock.Substring($segEndIdx)

$newXml = $head + $newMathBlock + $tail

$rng.InsertXML($newXml)

Write-Host "Equation patched."
